$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Label" column header - match the style of the other header cells
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Label values: 0 = Control rows, 1 = MDD rows (two blocks of 10 rows each)
$labels = @(0,0,0,0,0,1,1,1,1,1,0,0,0,0,0,1,1,1,1,1)
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $labels[$i]
}

# Refit values that changed slightly between the two runs
$ws.Range("D2").Value = 0.4892605029372821
$ws.Range("E2").Value = 0.4892605029372821

$ws.Range("D5").Value = 0.2065786190932306
$ws.Range("E5").Value = 0.2065786190932306

$ws.Range("D6").Value = 0.3505356885848804
$ws.Range("E6").Value = 0.3505356885848804

$ws.Range("D9").Value = 0.3132377843154457
$ws.Range("E9").Value = 0.6867622156845543

$ws.Range("D11").Value = 0.4561264791309875
$ws.Range("E11").Value = 0.5438735208690125
$ws.Range("F11").Value = 0.670553982257843
